$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New age-group headers in D1:P1
$headers = @(
    "Agrupacion de Edades 0-4",
    "Agrupacion de Edades 5-9",
    "Agrupacion de Edades 10-14",
    "Agrupacion de Edades 15-19",
    "Agrupacion de Edades 20-24",
    "Agrupacion de Edades 25-34",
    "Agrupacion de Edades 35-44",
    "Agrupacion de Edades 45-54",
    "Agrupacion de Edades 55-59",
    "Agrupacion de Edades 60-64",
    "Agrupacion de Edades 65-74",
    "Agrupacion de Edades 75-84",
    "Agrupacion de Edades 85+"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 4 + $i  # D = 4
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Formula in D2: sum of B2:B6 (calories for ages 0-4)
$ws.Range("D2").Formula = "=SUM(B2:B6)"

# Update selection to D7 (matches target sheetView)
$ws.Range("D7").Select()
